$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# The last 35 comments on Sheet1 (rows 276-310) are being relocated to the top
# of Sheet2 (new rows 1-35). Read them out first, in order, before they move.
$firstRow = 276
$lastRow  = 310
$rowCount = $lastRow - $firstRow + 1

$movedValues = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $movedValues += , $ws1.Range("A$r").Value()
}

# Drop those rows from the bottom of Sheet1 (dimension A1:A310 -> A1:A275).
$ws1.Rows("$($firstRow):$($lastRow)").Delete()

# Open up space at the top of Sheet2 for the relocated comments
# (dimension A1:A50 -> A1:A85).
$ws2.Rows("1:$rowCount").Insert()

# Write the relocated comments into the new rows, preserving their order.
for ($i = 0; $i -lt $movedValues.Count; $i++) {
    $destRow = $i + 1
    $ws2.Range("A$destRow").Value = $movedValues[$i]
}
